$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31 (old row31 "Total" -> 32, old row32 "footer" -> 33)
$ws.Rows.Item(31).Insert()

# Duplicate row 30's full formatting (styles + merged-cell layout) into the new row 31
# by copying it cell-for-cell; we'll overwrite the copied values right after.
$ws.Range("A30:Q30").Copy($ws.Range("A31:Q31"))
$ws.Rows.Item(31).RowHeight = 25.5

# Fill in the new sale-row content
$ws.Range("A31").Value = 25
$ws.Range("C31").Value = "مناديل سولو سحب صغيره"
$ws.Range("H31").Value = "1:0"
# L31 already reads "0" (copied verbatim from L30) - leave it untouched so the
# text stays a shared string instead of being reinterpreted as a number.
$ws.Range("N31").Value = "35.00"
# P31 must hold the text "35.0000" (not the number 35); assigning a numeric-looking
# string via .Value would get reinterpreted as a number, so copy an existing cell
# that already contains this exact text/style combination instead.
$ws.Range("P13").Copy($ws.Range("P31"))
$ws.Range("Q31").Value = "1:0"

# Update the Total row (now row 32) with the new sum
$ws.Range("P32").Value = 685.44000000000005

# Update the footer timestamp (now row 33) to the new generation time
$ws.Range("A33").Value = "Thursday, 7 August, 2025 12:10 PM"
